# Update cryptocurrency price/volume data per the latest symbol list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.488"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.31%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07984"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.988"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.08%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.82%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9499"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.91%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.19%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1881"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.10%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "10.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "25.75%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09976"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.15%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04827"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1063"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.16%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001280"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.76%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04087"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.57%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005936"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.97%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.365"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.99%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3468"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.99%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.53%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2588"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.54%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001266"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.00%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004265"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.23%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.80%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.05%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-2.47%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05646"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.57%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007547"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.73%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.00%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007422"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002017"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.31%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.53%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007140"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.54%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.10%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003533"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "55.66%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003795"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.76%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.10%"
